# Applies a weekly re-shuffle of the "Damasco" price rows: the date (D) and
# variety/quality/volume/price/unit/origin columns (K-T) get reassigned
# across rows 2,4,5,7,8,9,10,11,12 while A,B,C,E,F,G,H,I,J stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row, taken directly from the target data (columns D, K:T)
$rowData = @{
    2  = @{ D = 44544; K = "Castle Brite"; L = "Segunda"; M = 160; N = 16000; O = 17000; P = 16500; Q = "`$/caja 15 kilos"; R = "Región de O'Higgins"; S = 1100; T = 15 }
    4  = @{ D = 44537; K = "Castle Brite"; L = "Primera"; M = 60;  N = 21000; O = 21500; P = 21250; Q = "`$/caja 15 kilos"; R = "Región de O'Higgins"; S = 1417; T = 15 }
    5  = @{ D = 44168; K = "Castle Brite"; L = "Primera"; M = 30;  N = 16000; O = 17000; P = 16500; Q = "`$/caja 16 kilos granel"; R = "Región de Coquimbo"; S = 1031; T = 16 }
    7  = @{ D = 44174; K = "Castle Brite"; L = "Primera"; M = 75;  N = 9000;  O = 10000; P = 9467;  Q = "`$/caja 10 kilos"; R = "Región de O'Higgins"; S = 947;  T = 10 }
    8  = @{ D = 44551; K = "Castle Brite"; L = "Primera"; M = 120; N = 15500; O = 16000; P = 15750; Q = "`$/caja 15 kilos"; R = "Región de O'Higgins"; S = 1050; T = 15 }
    9  = @{ D = 44552; K = "Castle Brite"; L = "Primera"; M = 120; N = 15500; O = 16000; P = 15750; Q = "`$/caja 15 kilos"; R = "Región de O'Higgins"; S = 1050; T = 15 }
    10 = @{ D = 44187; K = "Dina";         L = "Primera"; M = 55;  N = 15000; O = 16000; P = 15455; Q = "`$/caja 15 kilos granel"; R = "Región de O'Higgins"; S = 1030; T = 15 }
    11 = @{ D = 44176; K = "Castle Brite"; L = "Primera"; M = 50;  N = 17000; O = 18000; P = 17400; Q = "`$/caja 18 kilos"; R = "Región de O'Higgins"; S = 967;  T = 18 }
    12 = @{ D = 44189; K = "Dina";         L = "Primera"; M = 80;  N = 16000; O = 17000; P = 16562; Q = "`$/caja 18 kilos"; R = "Región de O'Higgins"; S = 920;  T = 18 }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
    $ws.Range("T$r").Value = $vals.T
}
